$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update A110 date/time value ---
$ws.Cells.Item(110, 1).Value = 45447.2916666667

# --- Add new row 111 ---
# Copy A110's formatting (date/time number format + font) onto A111 so the
# new date cell carries the same style index as the rest of column A, then
# write its actual value.
$a111 = $ws.Cells.Item(111, 1)
$ws.Cells.Item(110, 1).Copy()
$a111.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$a111.Value = 45448.6494560185

$ws.Cells.Item(111, 2).Value = 68000
$ws.Cells.Item(111, 3).Value = 2.8199999332428
$ws.Cells.Item(111, 4).Value = 2.52999997138977
$ws.Cells.Item(111, 5).Value = 2.74000000953674
$ws.Cells.Item(111, 6).Value = 2.53999996185303

# adj_close (G111) must be stored as text "2.53999996185303", matching the
# source data's convention of quoting numeric-looking adj_close values.
$g111 = $ws.Cells.Item(111, 7)
$g111.NumberFormat = "@"
$g111.Value = "2.53999996185303"
$g111.Style = "Normal"

# ticker (H111) - text
$h111 = $ws.Cells.Item(111, 8)
$h111.Value = "LS.MI"
